$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header value tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 values were deleted entirely, C2 and E2 get new values
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -7.3733077377562868
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -3.0402514088270851

# Row 3 value tweaks
$ws.Range("B3").Value = -10.570762868888698
$ws.Range("C3").Value = 6.8726977810813397
$ws.Range("D3").Value = -4.6204314984218495
$ws.Range("E3").Value = 22.851515410309517

# Update the selection to match the new authored state
$ws.Range("B1:E3").Select()
